# Auto update: 2025-11-29 18:42:09
# Refresh the 국장_반도체_분석 sheet with the latest scored values.
# DB HiTek moves up into the #2 slot (right after Samsung Electronics),
# and the previous #2/#3 rows (240810.KS / 058470.KS) shift down to #3/#4.
# A new "5일수익률" (E column) figure is now populated for every ticker,
# and each row's rule score / probabilities / final score get refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : SamsungElec / 005930.KS  (stays in place, values refreshed)
$ws.Cells.Item(2, 2).Value = "SamsungElec"
$ws.Cells.Item(2, 3).Value = "005930.KS"
$ws.Cells.Item(2, 4).Value = 100500
$ws.Cells.Item(2, 5).Value = 49.9
$ws.Cells.Item(2, 6).Value = 6.01
$ws.Cells.Item(2, 7).Value = 30
$ws.Cells.Item(2, 8).Value = 63
$ws.Cells.Item(2, 9).Value = 73
$ws.Cells.Item(2, 10).Value = 73
$ws.Cells.Item(2, 11).Value = 64
$ws.Cells.Item(2, 13).Value = "📈 매수 관찰 구간입니다."
$ws.Cells.Item(2, 14).Value = 85.87127175646313

# Row 3 : DB HiTek / 000990.KS  (moved up from row 5)
$ws.Cells.Item(3, 2).Value = "DB HiTek"
$ws.Cells.Item(3, 3).Value = "000990.KS"
$ws.Cells.Item(3, 4).Value = 63600
$ws.Cells.Item(3, 5).Value = 53
$ws.Cells.Item(3, 6).Value = 2.75
$ws.Cells.Item(3, 7).Value = 40
$ws.Cells.Item(3, 8).Value = 36
$ws.Cells.Item(3, 9).Value = 56
$ws.Cells.Item(3, 10).Value = 50
$ws.Cells.Item(3, 11).Value = 60.2
$ws.Cells.Item(3, 13).Value = "📈 매수 관찰 구간입니다."
$ws.Cells.Item(3, 14).Value = 85.87127175646313

# Row 4 : 240810.KS  (moved down from row 3)
$ws.Cells.Item(4, 2).Value = "240810.KS,0P00017YB3,330568"
$ws.Cells.Item(4, 3).Value = "240810.KS"
$ws.Cells.Item(4, 4).Value = 61300
$ws.Cells.Item(4, 5).Value = 30.8
$ws.Cells.Item(4, 6).Value = 8.109999999999999
$ws.Cells.Item(4, 7).Value = 20
$ws.Cells.Item(4, 8).Value = 63
$ws.Cells.Item(4, 9).Value = 70
$ws.Cells.Item(4, 10).Value = 66
$ws.Cells.Item(4, 11).Value = 59.8
$ws.Cells.Item(4, 13).Value = "⛔ 관망하십시오."
$ws.Cells.Item(4, 14).Value = 85.87127175646313

# Row 5 : 058470.KS  (moved down from row 4)
$ws.Cells.Item(5, 2).Value = "058470.KS,0P0000ASU1,98886"
$ws.Cells.Item(5, 3).Value = "058470.KS"
$ws.Cells.Item(5, 4).Value = 68300
$ws.Cells.Item(5, 5).Value = 71.40000000000001
$ws.Cells.Item(5, 6).Value = 25.55
$ws.Cells.Item(5, 7).Value = 40
$ws.Cells.Item(5, 8).Value = 36
$ws.Cells.Item(5, 9).Value = 50
$ws.Cells.Item(5, 10).Value = 63
$ws.Cells.Item(5, 11).Value = 57.8
$ws.Cells.Item(5, 13).Value = "⛔ 관망하십시오."
$ws.Cells.Item(5, 14).Value = 85.87127175646313

# Row 6 : SK hynix / 000660.KS  (stays in place, values refreshed)
$ws.Cells.Item(6, 2).Value = "SK hynix"
$ws.Cells.Item(6, 3).Value = "000660.KS"
$ws.Cells.Item(6, 4).Value = 530000
$ws.Cells.Item(6, 5).Value = 35.6
$ws.Cells.Item(6, 6).Value = 1.8
$ws.Cells.Item(6, 7).Value = 20
$ws.Cells.Item(6, 8).Value = 60
$ws.Cells.Item(6, 9).Value = 60
$ws.Cells.Item(6, 10).Value = 70
$ws.Cells.Item(6, 11).Value = 55.8
$ws.Cells.Item(6, 13).Value = "⛔ 관망하십시오."
$ws.Cells.Item(6, 14).Value = 85.87127175646313

# Row 7 : 403870.KS  (stays in place, values refreshed)
$ws.Cells.Item(7, 2).Value = "403870.KS,0P0001PE9K,566428"
$ws.Cells.Item(7, 3).Value = "403870.KS"
$ws.Cells.Item(7, 4).Value = 30250
$ws.Cells.Item(7, 5).Value = 39.7
$ws.Cells.Item(7, 6).Value = 6.7
$ws.Cells.Item(7, 7).Value = 20
$ws.Cells.Item(7, 8).Value = 53
$ws.Cells.Item(7, 9).Value = 46
$ws.Cells.Item(7, 10).Value = 56
$ws.Cells.Item(7, 11).Value = 50.2
$ws.Cells.Item(7, 13).Value = "⛔ 관망하십시오."
$ws.Cells.Item(7, 14).Value = 85.87127175646313
